$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value  = -12.71440000000001
$ws.Range("A8").Value  = -21.1868
$ws.Range("A10").Value = -20.48349999999998
$ws.Range("A12").Value = -22.28220000000004
$ws.Range("B13").Value = 6.281899999999997
$ws.Range("A18").Value = -22.33910000000003
$ws.Range("C20").Value = -13.40359999999999
$ws.Range("A25").Value = -22.24940000000003
